$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 0.375
$ws.Range("C1").Value = 0.01602287055634733
$ws.Range("D1").Value = 0
$ws.Range("E1").Value = 162.3018583932938

$ws.Range("B2").Value = 0.5
$ws.Range("C2").Value = 0.007981730387757634
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 175.2116728691194
